$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AU-4,AU-4 (1)"
$ws.Range("A3").Value = "AU-4,AU-14 (1)"
$ws.Range("A4").Value = "AU-4,CM-6 b"
$ws.Range("A5").Value = "SC-5,CM-6 b,SC-5 (2)"
$ws.Range("A6").Value = "AU-12 (3),AU-8 b,CM-5 (1),AC-6 (9),AC-6 (8),AU-7 a,AU-7 b"
$ws.Range("A7").Value = "AU-12 (3),AU-12 c,AU-8 b,AU-12 a,CM-5 (1),AU-7 a,CM-6 b,AU-7 b"
$ws.Range("A8").Value = "CM-6 b,AC-17 (1),CM-7 b,AC-17 (9)"
$ws.Range("A14").Value = "CM-7 (2),CM-7 (5) (b)"
$ws.Range("A15").Value = "CM-7 (2),CM-7 (5) (b)"
$ws.Range("A17").Value = "CM-7 (2),CM-6 b"
$ws.Range("A22").Value = "CM-7 (2),CM-6 b"
$ws.Range("A23").Value = "CM-7 (2),CM-6 b"
$ws.Range("A38").Value = "AC-7 b,AC-7 a"
$ws.Range("A39").Value = "AC-7 b,AC-7 a"
$ws.Range("A40").Value = "AC-7 b,AC-7 a"
$ws.Range("A41").Value = "AC-7 b,AC-7 a"
$ws.Range("A45").Value = "AU-3 (1),IA-8,IA-2"
$ws.Range("A46").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A47").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A48").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A49").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A50").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A51").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A52").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A53").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A54").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A55").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A56").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A57").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A58").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A59").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A60").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A61").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A62").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A63").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A64").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A65").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A66").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A67").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A68").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A69").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A70").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A71").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A72").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A73").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A74").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A75").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A76").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A77").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A78").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A79").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A80").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A81").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A82").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A83").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A84").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A85").Value = "AU-3,MA-4 (1) (a),AU-3 (1)"
$ws.Range("A86").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A87").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A88").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A89").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A90").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A91").Value = "AU-12 c,MA-4 (1) (a),AU-3 (1)"
$ws.Range("A92").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A93").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A94").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A95").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A96").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A97").Value = "AU-3,AU-12 c,MA-4 (1) (a),AU-3 (1)"
$ws.Range("A98").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A99").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A100").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A101").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A102").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A103").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A104").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A105").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A106").Value = "AC-2 (4),AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A107").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A108").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A109").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A110").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A111").Value = "AC-2 (4),AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3"
$ws.Range("A112").Value = "AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-14 (1),AU-3"
$ws.Range("A115").Value = "AC-11 b,AC-6 (10)"
$ws.Range("A120").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c"
$ws.Range("A121").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c"
$ws.Range("A122").Value = "AU-9,AU-12 c"
$ws.Range("A127").Value = "AC-6 (9),AC-2 (4),CM-5 (1),AU-12 c"
$ws.Range("A129").Value = "IA-5 (1) (a),CM-6 b,IA-5 (1) (b)"
$ws.Range("A133").Value = "SC-8,AC-17 (2),SC-13,MA-4 c"
$ws.Range("A135").Value = "SC-10,AC-12"
$ws.Range("A136").Value = "SC-10,AC-12"
$ws.Range("A137").Value = "AC-11 a,SC-10"
$ws.Range("A138").Value = "AU-12 a,AU-6 (4),AU-3 (1),CM-5 (1),AU-7 (1),AU-7 a,CM-6 b,MA-4 (1) (a),AU-14 (1),AU-3"
$ws.Range("A141").Value = "AU-9 (3),AU-9"
$ws.Range("A142").Value = "AU-9 (3),AU-9"
$ws.Range("A143").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A144").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A145").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A146").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A147").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A148").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A149").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A150").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A151").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A152").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A153").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A154").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A155").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A166").Value = "SC-8 (2),SC-8,SC-8 (1)"
$ws.Range("A168").Value = "SC-8,AC-17 (2)"
$ws.Range("A171").Value = "AC-6 (9),AC-2 (4),AU-12 c"
$ws.Range("A173").Value = "AC-11 a,AC-11 b"
$ws.Range("A180").Value = "AU-4 (1),CM-6 b,AU-6 (4)"
$ws.Range("A181").Value = "CM-6 b,AC-17 (1),CM-7 b"
$ws.Range("A194").Value = "AU-3,CM-6 b"
$ws.Range("A200").Value = "AU-3,AU-4 (1)"
$ws.Range("A207").Value = "AU-4 (1),CM-6 b"
$ws.Range("A208").Value = "SC-28 (1),SC-28"
$ws.Range("A220").Value = "IA-2 (5),CM-6 b"
$ws.Range("A221").Value = "IA-2 (4),IA-2,IA-2 (2),IA-2 (5),IA-2 (3)"
$ws.Range("A222").Value = "IA-2 (4),IA-2,IA-2 (2),IA-2 (5),IA-2 (3)"
$ws.Range("A223").Value = "AC-18 (1),SC-8,SC-8 (1)"
$ws.Range("A225").Value = "IA-5 (1) (c),IA-7"
$ws.Range("A230").Value = "IA-7,CM-7 a"
$ws.Range("A244").Value = "SI-16,CM-6 b,SC-2"
$ws.Range("A269").Value = "IA-2 (2),CM-6 b"
$ws.Range("A276").Value = "SC-2,SC-4"
$ws.Range("A277").Value = "SC-2,SC-4"
$ws.Range("A299").Value = "IA-2 (11),IA-2 (1),IA-2 (12)"
$ws.Range("A309").Value = "AU-8 (1) (b),AU-8 b,AU-8 (1) (a)"
$ws.Range("A329").Value = "AU-5 a,AU-5 b"
$ws.Range("A341").Value = "CM-7 b,IA-3"
$ws.Range("A342").Value = "CM-7 b,CM-7 a"
$ws.Range("A343").Value = "CM-7 b,CM-7 a"
$ws.Range("A345").Value = "AC-18 (1),CM-7 a"
$ws.Range("A346").Value = "IA-5 (1) (c),CM-6 b,CM-7 a"
$ws.Range("A357").Value = "AC-11 b,AC-11 (1)"
$ws.Range("A360").Value = "SI-6 b,SI-6 d,CM-3 (5)"
$ws.Range("A366").Value = "SI-16,CM-7 a"
$ws.Range("A373").Value = "CM-6 b,CM-7 a"
$ws.Range("A374").Value = "CM-6 b,CM-7 a"
$ws.Range("A375").Value = "CM-6 b,CM-7 a"
$ws.Range("A388").Value = "SI-6 a,SC-3"
$ws.Range("A390").Value = "IA-5 (1) (a),CM-6 b"
$ws.Range("A396").Value = "SI-6 d,CM-3 (5)"
$ws.Range("A397").Value = "SI-16,CM-6 b"
$ws.Range("A447").Value = "IA-5 (1) (c),CM-6 b"
